$d = $word.ActiveDocument

# Step 1: Delete the "_GoBack" bookmark from its old position
# (right after "...exploration policy (more on this later.)")
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# Step 2: Replace the two-run italic "Demo" + "s" with a single italic run "Agents"
$rngReplace = $d.Content
$rngReplace.Find.Execute("Demos", $true, $false, $false, $false, $false, $true, 1, $false, "Agents", 2)

# Step 3: Insert the bookmark back, now right after the "Agents" run (before " directory.")
$findDir = $d.Content
$findDir.Find.Execute(" directory.", $false)
$bmPoint = $d.Range($findDir.Start, $findDir.Start)
$d.Bookmarks.Add("_GoBack", $bmPoint)

Write-Output "done"
